$d = $word.ActiveDocument

# 1) Merge the first paragraph's two runs (and remove the _GoBack bookmark in between)
#    into a single run of text "PROBANDO GIT LOCAL, DESDE ESTE ARCHIVO."
$d.Content.Find.Execute("PROBANDO GIT LOCAL, DESDE ESTE ARCHIVO.", $false, $false, $false, $false, $false, $true, 1, $false, "PROBANDO GIT LOCAL, DESDE ESTE ARCHIVO.", 2)

# 2) Position at the end of that first paragraph and insert two new paragraphs:
#    one empty, and one with "Agregando cambios por Alejandro"
$p1 = $d.Paragraphs.Item(1)
$endRange = $p1.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertBefore("Agregando cambios por Alejandro")
